# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Price" column cells are plain text (e.g. "42.672.36", "254.09"), not numbers.
# Flip to the Text format before writing so Excel does not auto-coerce a
# numeric-looking string into a Number cell, then drop back to the Normal
# style so no stray formatting is left behind.
$priceCells = @("D2", "D3", "D5", "D7", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D22", "D24", "D26", "D27", "D31", "D36", "D37", "D38", "D39", "D41", "D43", "D45", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.672.36'
$ws.Range("E2").Value = '  -0.82%  '

$ws.Range("D3").Value = '2.208.72'
$ws.Range("E3").Value = '  -1.40%  '

$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").Value = '254.09'
$ws.Range("E5").Value = '  +3.31%  '

$ws.Range("D7").Value = '74.94'
$ws.Range("E7").Value = '  -1.82%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  -4.61%  '

$ws.Range("D10").Value = '40.93'
$ws.Range("E10").Value = '  -0.02%  '

$ws.Range("D11").Value = '0.0917'
$ws.Range("E11").Value = '  -1.83%  '

$ws.Range("E12").Value = '  -1.28%  '

$ws.Range("D13").Value = '0.101'
$ws.Range("E13").Value = '  -0.44%  '

$ws.Range("D14").Value = '2.539.79'
$ws.Range("E14").Value = '  -0.55%  '

$ws.Range("D15").Value = '14.28'
$ws.Range("E15").Value = '  -2.53%  '

$ws.Range("D16").Value = '2.212.49'
$ws.Range("E16").Value = '  -0.85%  '

$ws.Range("D17").Value = '0.779'
$ws.Range("E17").Value = '  -4.11%  '

$ws.Range("D18").Value = '42.597.16'
$ws.Range("E18").Value = '  -0.78%  '

$ws.Range("E19").Value = '  -2.17%  '

$ws.Range("D20").Value = '71.02'
$ws.Range("E20").Value = '  -0.33%  '

$ws.Range("E21").Value = '  -1.01%  '

$ws.Range("D22").Value = '227.01'
$ws.Range("E22").Value = '  -1.01%  '

$ws.Range("E23").Value = '  -3.24%  '

$ws.Range("D24").Value = '9.38'

$ws.Range("D26").Value = '10.56'
$ws.Range("E26").Value = '  -3.09%  '

$ws.Range("D27").Value = '39.61'
$ws.Range("E27").Value = '  +3.42%  '

$ws.Range("E28").Value = '  +0.27%  '

$ws.Range("E29").Value = '  +2.84%  '

$ws.Range("E30").Value = '  -3.24%  '

$ws.Range("D31").Value = '172.94'
$ws.Range("E31").Value = '  -0.73%  '

$ws.Range("E32").Value = '  -0.71%  '

$ws.Range("E33").Value = '  +4.31%  '

$ws.Range("E34").Value = '  -3.23%  '

$ws.Range("E35").Value = '  -1.39%  '

$ws.Range("D36").Value = '0.108'
$ws.Range("E36").Value = '  -4.86%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '4.28'
$ws.Range("E37").Value = '  -2.55%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.0341'
$ws.Range("E38").Value = '  +4.37%  '

$ws.Range("D39").Value = '12.33'
$ws.Range("E39").Value = '  -5.27%  '

$ws.Range("E40").Value = '  -2.50%  '

$ws.Range("D41").Value = '2.74'
$ws.Range("E41").Value = '  +18.39%  '

$ws.Range("E42").Value = '  -5.84%  '

$ws.Range("D43").Value = '59.79'
$ws.Range("E43").Value = '  -0.46%  '

$ws.Range("E44").Value = '  -4.06%  '

$ws.Range("D45").Value = '101.38'
$ws.Range("E45").Value = '  -4.08%  '

$ws.Range("E46").Value = '  -3.55%  '

$ws.Range("E47").Value = '  -1.29%  '

$ws.Range("E48").Value = '  +3.92%  '

$ws.Range("E49").Value = '  -0.18%  '

$ws.Range("E50").Value = '  -1.26%  '

$ws.Range("D51").Value = '2.434.78'
$ws.Range("E51").Value = '  -0.25%  '

# Restore the Normal style on the price cells now that the text value is set.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}